# Generate Report for Handback
# Updates the localization-status workbook after a handback completes:
#  - Overview sheet: status text for a.md / b.md -> "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows, and widen columns that now hold
#    longer text.

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"
$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c48e3eba6053ec07966460dcf971dea04e84d329/e2e"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both file rows
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusDone
$wsOverview.Range("F2").Value = $statusDone
$wsOverview.Range("E3").Value = $statusDone
$wsOverview.Range("F3").Value = $statusDone

# Widen the now-longer status columns to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet (rows 2 & 3 = a.md / b.md)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusDone
$wsZh.Range("C3").Value = $statusDone

# Latest Handback DateTime (column K) gets the real handback timestamp.
$wsZh.Range("K2").Value = "2016-08-19 06:36:27"
$wsZh.Range("K3").Value = "2016-08-19 06:36:27"

# Latest Handback File (column J) gets the generated xliff file name.
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Latest Target File (column I) now links to the localized a.md, same on both rows.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$githubBase/a.md", "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$githubBase/a.md", "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$githubBase/b.md", "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$githubBase/a.md", "", "", "a.md")

# Column widths: Status (C) and Latest Handback File (J) need more room now.
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet (rows 2 & 3 = a.md / b.md)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusDone
$wsDe.Range("C3").Value = $statusDone

# Latest Handback DateTime (column K) gets the real handback timestamp.
$wsDe.Range("K2").Value = "2016-08-19 06:36:34"
$wsDe.Range("K3").Value = "2016-08-19 06:36:34"

# Latest Handback File (column J) gets the generated xliff file name.
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# Latest Target File (column I) now links to the localized a.md, same on both rows.
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$githubBase/a.md", "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$githubBase/a.md", "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$githubBase/b.md", "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$githubBase/a.md", "", "", "a.md")

# Column widths: Status (C) and Latest Handback File (J) need more room now.
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(10).ColumnWidth = 40
